# Unit5_zodiac.docx edit: the body copy was re-pasted as clean OOXML,
# which drops the per-run/per-paragraph direct formatting (w:pPr / w:rPr
# with the Helvetica/2D3B45 styling) that every paragraph had picked up
# from the original source. The wording is unchanged except that the
# last sentence of paragraph 2 now has the cursor (the "_GoBack" bookmark)
# parked after "...database operatio|ns." and the trailing empty
# paragraph loses that bookmark.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# One <w:p> fragment per paragraph, stripped of all w:pPr/w:rPr - only
# the raw runs (plus the spell-check markers and the relocated bookmark)
# survive.
$paragraphsXml = @(
  "<w:p xmlns:w='$wNs'><w:r><w:t xml:space='preserve'>Two SQL files have been included in the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>src</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>/main/resources folder. Using the schema provided, modify the application so that it reads Fortunes and Months from the database instead of from JSON files.</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t xml:space='preserve'>Modify the entity classes to include persistence annotations. Create </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>CrudRepository</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> classes that the Service classes will use for database operatio</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/><w:r><w:t>ns.</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t>Submit the URL of your GitHub repository.</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t>$([char]0x00A0)</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t>Notes</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t xml:space='preserve'>The project has already been configured for you (Maven dependencies and </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>application.properties</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>)</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'><w:r><w:t>You do not have to perform validation.</w:t></w:r></w:p>",

  "<w:p xmlns:w='$wNs'/>"
)

$bodyXml = [string]::Join("", $paragraphsXml)

# Content spans every paragraph up to (but not including) the sectPr,
# so replacing it in one shot swaps the whole body while leaving the
# section properties (page size/margins/etc.) untouched.
$d.Content.InsertXML($bodyXml)
